$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit inserts one new data row (a new weekly observation) right after
# the existing row 42, pushing the previous rows 43..153 down to 44..154.
$ws.Rows(43).Insert()

# Populate the newly inserted row 43 with the new observation's data.
$ws.Cells.Item(43, 1).Value  = 7
$ws.Cells.Item(43, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(43, 3).Value  = "Ñuble"
$ws.Cells.Item(43, 4).Value  = 44526
$ws.Cells.Item(43, 5).Value  = 16
$ws.Cells.Item(43, 6).Value  = 100112003
$ws.Cells.Item(43, 7).Value  = "Ajo"
$ws.Cells.Item(43, 8).Value  = "Chino"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 60
$ws.Cells.Item(43, 11).Value = 17000
$ws.Cells.Item(43, 12).Value = 18000
$ws.Cells.Item(43, 13).Value = 17500
$ws.Cells.Item(43, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(43, 15).Value = "China"
$ws.Cells.Item(43, 16).Value = 1750
$ws.Cells.Item(43, 17).Value = 10
$ws.Cells.Item(43, 18).Value = "Hortaliza"
